# Updates Kane Williamson's per-innings batting activity (runs/balls/fours/sixes)
# in rows 2-12 to reflect the latest figures pulled from the Excel form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (runs, balls, fours, sixes)
$data = @{
  2  = @("50","44","2","2")
  3  = @("8","14","0","0")
  4  = @("11","10","1","0")
  5  = @("67","45","5","4")
  6  = @("57","39","7","0")
  7  = @("29","19","4","1")
  8  = @("22","12","0","2")
  9  = @("9","13","1","0")
  10 = @("3","5","0","0")
  11 = @("41","26","5","0")
  12 = @("20","10","1","1")
}

$cols = @("C","D","E","F")

# The source values are text (numbers stored as text), so force the range to
# Text format before writing, otherwise Excel would auto-convert the numeric
# looking strings into real numbers.
$rng = $ws.Range("C2:F12")
$rng.NumberFormat = "@"

for ($r = 2; $r -le 12; $r++) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}

# Drop the temporary Text number format again so the cell styling matches the
# original workbook (only the values themselves changed).
$rng.ClearFormats()
